$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from H1 to I1:J1 before assigning values
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for I2:J71
$data = @(
    @(2, 8, 8),
    @(3, 6, 6),
    @(4, 6, 6),
    @(5, 5, 5),
    @(6, 6, 7),
    @(7, 8, 8),
    @(8, 7, 7),
    @(9, 7, 7),
    @(10, 7, 7),
    @(11, 8, 8),
    @(12, 7, 7),
    @(13, 6, 6),
    @(14, 6, 6),
    @(15, 9, 9),
    @(16, 9, 9),
    @(17, 8, 8),
    @(18, 8, 8),
    @(19, 7, 7),
    @(20, 7, 8),
    @(21, 8, 8),
    @(22, 9, 9),
    @(23, 6, 6),
    @(24, 7, 7),
    @(25, 6, 6),
    @(26, 7, 7),
    @(27, 8, 8),
    @(28, 7, 7),
    @(29, 6, 6),
    @(30, 8, 8),
    @(31, 8, 8),
    @(32, 6, 6),
    @(33, 7, 7),
    @(34, 7, 7),
    @(35, 8, 8),
    @(36, 9, 9),
    @(37, 8, 8),
    @(38, 8, 8),
    @(39, 7, 7),
    @(40, 7, 8),
    @(41, 6, 6),
    @(42, 7, 7),
    @(43, 8, 8),
    @(44, 7, 7),
    @(45, 6, 7),
    @(46, 7, 8),
    @(47, 9, 9),
    @(48, 8, 8),
    @(49, 8, 8),
    @(50, 7, 7),
    @(51, 7, 8),
    @(52, 8, 8),
    @(53, 7, 7),
    @(54, 7, 7),
    @(55, 8, 8),
    @(56, 6, 6),
    @(57, 8, 8),
    @(58, 7, 7),
    @(59, 8, 8),
    @(60, 8, 8),
    @(61, 7, 7),
    @(62, 7, 7),
    @(63, 8, 8),
    @(64, 6, 6),
    @(65, 7, 7),
    @(66, 7, 7),
    @(67, 7, 7),
    @(68, 8, 9),
    @(69, 8, 8),
    @(70, 6, 6),
    @(71, 3, 3)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 9).Value = $entry[1]
    $ws.Cells.Item($r, 10).Value = $entry[2]
}
